$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 171-172, pushing the existing rows 171-264 down to 173-266.
$ws.Rows("171:172").Insert()

# New row 171: Primera, new later date, new volume/price data
$ws.Range("A171").Value = 4
$ws.Range("B171").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C171").Value = "Los Lagos"
$ws.Range("D171").Value = 44572
$ws.Range("E171").Value = 10
$ws.Range("F171").Value = 100112023
$ws.Range("G171").Value = "Brócoli"
$ws.Range("H171").Value = "Sin especificar"
$ws.Range("I171").Value = "Primera"
$ws.Range("J171").Value = 400
$ws.Range("K171").Value = 1400
$ws.Range("L171").Value = 1400
$ws.Range("M171").Value = 1400
$ws.Range("N171").Value = "`$/unidad"
$ws.Range("O171").Value = "Región Metropolitana"
$ws.Range("P171").Value = 1400
$ws.Range("Q171").Value = 1
$ws.Range("R171").Value = "Hortaliza"

# New row 172: Segunda, same new date
$ws.Range("A172").Value = 4
$ws.Range("B172").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C172").Value = "Los Lagos"
$ws.Range("D172").Value = 44572
$ws.Range("E172").Value = 10
$ws.Range("F172").Value = 100112023
$ws.Range("G172").Value = "Brócoli"
$ws.Range("H172").Value = "Sin especificar"
$ws.Range("I172").Value = "Segunda"
$ws.Range("J172").Value = 400
$ws.Range("K172").Value = 1000
$ws.Range("L172").Value = 1000
$ws.Range("M172").Value = 1000
$ws.Range("N172").Value = "`$/unidad"
$ws.Range("O172").Value = "Región Metropolitana"
$ws.Range("P172").Value = 1000
$ws.Range("Q172").Value = 1
$ws.Range("R172").Value = "Hortaliza"
